# Saldo.xlsx update: refresh the account balance export.
#
# Net effect vs. the previous export:
#   - Removed accounts 005055226 (DANILO) and 004459875 (HELVECIO) - no longer present.
#   - Removed account 005995120 (ERIK) - no longer present.
#   - Account 004511696 (KRYSCIA) balance updated from 6610.91 to 610.91, and the
#     row now sits further down the list (after ROBERTO / 004207658, before LUZIMAR).
#
# All other rows are unchanged. We apply this as: delete the three obsolete rows
# (bottom-up, by their original row numbers, so earlier deletes don't shift the
# row numbers of rows still to be removed), then insert a fresh row for KRYSCIA
# at its new location with the updated balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so previously-computed row numbers stay valid.
$ws.Rows.Item(19).Delete()              # 005995120 / ERIK / 3090.73
$ws.Rows.Item(15).Delete()              # 004511696 / KRYSCIA / 6610.91 (old position/value)
$ws.Range("A12:A13").EntireRow.Delete() # 005055226 / DANILO / 13227.6 ; 004459875 / HELVECIO / 11028.78

# Insert the refreshed KRYSCIA row right after ROBERTO (row 20) and before LUZIMAR (row 21).
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "004511696"
$ws.Cells.Item(21, 2).Value = "KRYSCIA"
$ws.Cells.Item(21, 3).Value = 610.91
